# Replace the single paragraph's content: the old text "12345" (with its
# language-tagged run properties and the _GoBack bookmark) is swapped out
# for the new multi-run "Teste do primeiro commit realizado em command line
# git" text, including the <w:proofErr> spell-check markers Word leaves
# around the terms it does not recognise.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1079A16D" w14:textId="380DD983" w:rsidR="004C42EC" w:rsidRPr="00E74D95" w:rsidRDefault="00E74D95"><w:r><w:t xml:space="preserve">Teste do primeiro </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> realizado em </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>co</w:t></w:r><w:r><w:t>mmand</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>line</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@

$r.InsertXML($xml)
